$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.251060366630554
$ws.Range("B1").Value = 1.765236139297485
$ws.Range("C1").Value = 1.91124427318573
$ws.Range("D1").Value = 7.081009387969971
$ws.Range("E1").Value = 1.670521020889282
